$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price and Volume columns to Text format to preserve exact string values
# (prevents Excel auto-converting numeric-looking strings to numbers)
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '23.319.95'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.633.03'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.48%  '

$ws.Range("D5").Value = '1.004'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").Value = '0.3812'
$ws.Range("E7").Value = '  +1.12%  '

$ws.Range("D8").Value = '51.82'
$ws.Range("E8").Value = '  -0.57%  '

$ws.Range("D9").Value = '0.3595'
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").Value = '0.08148'
$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("D11").Value = '1.219'
$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").Value = '22.34'
$ws.Range("E13").Value = '  -1.53%  '

$ws.Range("D14").Value = '6.399'
$ws.Range("E14").Value = '  -2.47%  '

$ws.Range("D15").Value = '7.293'
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D16").Value = '0.00001229'
$ws.Range("E16").Value = '  -1.11%  '

$ws.Range("D17").Value = '1.634.84'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '94.90'
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").Value = '0.06951'
$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("D20").Value = '6.573'
$ws.Range("E20").Value = '  +1.73%  '

$ws.Range("D21").Value = '17.33'
$ws.Range("E21").Value = '  -3.47%  '

$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").Value = '12.43'
$ws.Range("E23").Value = '  -2.30%  '

$ws.Range("D24").Value = '23.412.22'
$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").Value = '2.542'
$ws.Range("E25").Value = '  +4.03%  '

$ws.Range("D26").Value = '3.034'
$ws.Range("E26").Value = '  -5.53%  '

$ws.Range("D27").Value = '21.17'
$ws.Range("E27").Value = '  +0.50%  '

$ws.Range("D28").Value = '151.45'
$ws.Range("E28").Value = '  +1.01%  '

$ws.Range("D29").Value = '5.271'
$ws.Range("E29").Value = '  -0.55%  '

$ws.Range("D30").Value = '133.72'
$ws.Range("E30").Value = '  -0.84%  '

$ws.Range("D31").Value = '1.818.15'
$ws.Range("E31").Value = '  +0.69%  '

$ws.Range("D32").Value = '1.084'
$ws.Range("E32").Value = '  +14.35%  '

$ws.Range("D33").Value = '2.148'
$ws.Range("E33").Value = '  -6.72%  '

$ws.Range("D34").Value = '6.456'
$ws.Range("E34").Value = '  -4.72%  '

$ws.Range("D35").Value = '11.34'
$ws.Range("E35").Value = '  +4.02%  '

$ws.Range("D36").Value = '0.02758'
$ws.Range("E36").Value = '  -2.31%  '

$ws.Range("D37").Value = '0.2500'
$ws.Range("E37").Value = '  -1.18%  '

$ws.Range("D38").Value = '0.08755'
$ws.Range("E38").Value = '  -0.79%  '

$ws.Range("D39").Value = '0.06990'
$ws.Range("E39").Value = '  -1.63%  '

$ws.Range("D40").Value = '5.911'
$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.346'
$ws.Range("E41").Value = '  -1.09%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.7003'
$ws.Range("E42").Value = '  -0.59%  '

$ws.Range("D43").Value = '12.11'
$ws.Range("E43").Value = '  -1.93%  '

$ws.Range("D44").Value = '15.58'
$ws.Range("E44").Value = '  -3.69%  '

$ws.Range("D45").Value = '0.6457'
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").Value = '1.002'
$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").Value = '2.280'
$ws.Range("E47").Value = '  -1.64%  '

$ws.Range("D48").Value = '3.960'
$ws.Range("E48").Value = '  -0.69%  '

$ws.Range("D49").Value = '0.07938'
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("D50").Value = '127.88'
$ws.Range("E50").Value = '  +1.28%  '

$ws.Range("D51").Value = '1.184'
$ws.Range("E51").Value = '  -1.84%  '

# Restore default style (remove the temporary text-format override)
$priceVolRange.Style = "Normal"